# Update the two customer-id strings wherever they occur on Sheet1,
# then update the selection on Sheet1 to B13:B14 with active cell B13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace cus_H3kv7w75JgdeGy -> cus_H4UfXFvcyIAryp (D3, D8, B13)
$ws.Range("D3").Value = "cus_H4UfXFvcyIAryp"
$ws.Range("D8").Value = "cus_H4UfXFvcyIAryp"
$ws.Range("B13").Value = "cus_H4UfXFvcyIAryp"

# Replace cus_H0So4YQfL0Mv1q -> cus_H4Uaiet9m8mye3 (D4, D9, B14)
$ws.Range("D4").Value = "cus_H4Uaiet9m8mye3"
$ws.Range("D9").Value = "cus_H4Uaiet9m8mye3"
$ws.Range("B14").Value = "cus_H4Uaiet9m8mye3"

# Update the selection / active cell shown in the sheet view.
$ws.Activate()
$ws.Range("B13:B14").Select()
$excel.ActiveCell = $ws.Range("B13")
